$d = $word.ActiveDocument

# 1) Insert a new list paragraph right after the final paragraph ("Finalizado el
#    turno de cada jugador,  se analizan los resultados finales."), then swap the
#    contents so the new sentence about player actions during the turn ends up
#    before the "Finalizado..." sentence - matching a normal "place cursor at the
#    end of 'Se analizan...', press Enter, type the new sentence" edit followed by
#    moving the trailing copy.
$pFinal = $d.Paragraphs(17)
$finalText = $pFinal.Range.Text
$pFinal.Range.InsertParagraphAfter()

$pFinalCopy = $d.Paragraphs(18)
$pFinalCopy.Range.Text = $finalText

$pNueva = $d.Paragraphs(17)
$pNueva.Range.Text = "Dependiendo de las condiciones en que se encuentra el jugador, durante la jugada el mismo podrá PEDIR CARTA, PLANTARSE, DOBLAR, SEPARAR."

# 2) Move the "_GoBack" bookmark to the end of the "Una vez finalizadas las
#    apuestas..." paragraph (excluding the trailing paragraph mark). Re-adding a
#    bookmark with the same name simply relocates it, regardless of where (if
#    anywhere) it currently sits.
$pApuestas = $d.Paragraphs(15)
$start = $pApuestas.Range.Start
$end = $pApuestas.Range.End
$bmRange = $d.Range($start, $end - 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
